$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing cell values
$ws.Range("A2").Value = "001q000000hmXhAAAU"
$ws.Range("A3").Value = "001q000000hmfgoAAA"
$ws.Range("C3").Value = "01tq0000001jhI0"
$ws.Range("D3").Value = "a0Nq0000003PBEa"
$ws.Range("K3").Value = "SVMXC__Service_Request__c SR_1 = new SVMXC__Service_Request__c(SVMXC__Status__c = 'Closed', SVMXC__Priority__c = 'High' );insert SR_1 ;"

# Match style (wrap text, default font) of other data cells for the newly populated cells
$ws.Range("C3").WrapText = $true
$ws.Range("D3").WrapText = $true

# Column J width change
$ws.Columns.Item(10).ColumnWidth = 20

# Row 3 height change
$ws.Rows.Item(3).RowHeight = 64

# Selection / view changes
$ws.Range("A3").Select()
